# Arreglos en el menu.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows ---
# Row 5 (D5): Price 102 -> 432
$ws.Cells.Item(5, 4).Value2 = 432

# Row 6 (C6): Author "Jhon Ricardo" -> "Jose Ramirez"
$ws.Cells.Item(6, 3).Value2 = "Jose Ramirez"

# --- Append new rows 12-19 (two repeated blocks of 4 books) ---
$newBooks = @(
    @("El que se duerme pierde", "Tom Peter", 16),
    @("Sin lugar a duda", "Ana Gutierrez", 26),
    @("El arte de dormir", "Nico", 32),
    @("Buscando a Nemo", "Humble Po", 41)
)

$row = 12
for ($block = 0; $block -lt 2; $block++) {
    foreach ($book in $newBooks) {
        $no = $row - 1
        $ws.Cells.Item($row, 1).Value2 = $no
        $ws.Cells.Item($row, 2).Value2 = $book[0]
        $ws.Cells.Item($row, 3).Value2 = $book[1]
        $ws.Cells.Item($row, 4).Value2 = $book[2]
        $row++
    }
}
